# Scheduled market-data refresh for the Leve profit sheets.
# Updates cached currentAveragePrice* / LevePrice* / LeveProfit* columns
# (H:N) per leve row across all eight crafting-job sheets with freshly
# pulled Universalis price data. Column layout per sheet:
#   H=currentAveragePrice  I=currentAveragePriceNQ  J=currentAveragePriceHQ
#   K=LevePriceNQ  L=LevePriceHQ  M=LeveProfitNQ  N=LeveProfitHQ

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6: Days of Chunder / Antidote
$ws.Cells.Item(6, 8).Value = 2166.6667  # H6: 3000 -> 2166.6667
$ws.Cells.Item(6, 9).Value = 500  # I6: 0 -> 500
$ws.Cells.Item(6, 11).Value = 1500  # K6: 0 -> 1500
$ws.Cells.Item(6, 13).Value = -1388  # M6: None -> -1388
# Row 11: Gotta Bounce / Rubber
$ws.Cells.Item(11, 8).Value = 68762.44  # H11: 57910.473 -> 68762.44
$ws.Cells.Item(11, 9).Value = 68762.44  # I11: 57910.473 -> 68762.44
$ws.Cells.Item(11, 11).Value = 68762.44  # K11: 57910.473 -> 68762.44
$ws.Cells.Item(11, 13).Value = -68622.44  # M11: -57770.473 -> -68622.44
# Row 17: One for the Road / Potion
$ws.Cells.Item(17, 8).Value = 1958.1904  # H17: 1930.4546 -> 1958.1904
$ws.Cells.Item(17, 10).Value = 1958.1904  # J17: 1930.4546 -> 1958.1904
$ws.Cells.Item(17, 12).Value = 5874.5712  # L17: 5791.3638 -> 5874.5712
$ws.Cells.Item(17, 14).Value = -6210.5712  # N17: -6127.3638 -> -6210.5712
# Row 28: The Writing Is Not on the Wall / Enchanted Silver Ink
$ws.Cells.Item(28, 8).Value = 1758.5555  # H28: 1763.7778 -> 1758.5555
$ws.Cells.Item(28, 10).Value = 2086.4167  # J28: 2094.25 -> 2086.4167
$ws.Cells.Item(28, 12).Value = 2086.4167  # L28: 2094.25 -> 2086.4167
$ws.Cells.Item(28, 14).Value = -3056.4167  # N28: -3064.25 -> -3056.4167
# Row 107: Another Man's Ink / Enchanted Truegold Ink
$ws.Cells.Item(107, 8).Value = 262.53845  # H107: 272 -> 262.53845
$ws.Cells.Item(107, 9).Value = 241.4  # I107: 251.66667 -> 241.4
$ws.Cells.Item(107, 11).Value = 241.4  # K107: 251.66667 -> 241.4
$ws.Cells.Item(107, 13).Value = 1678.6  # M107: 1668.33333 -> 1678.6
# Row 116: Growing Up / Growth Formula Kappa
$ws.Cells.Item(116, 8).Value = 37208548  # H116: 34728200 -> 37208548
$ws.Cells.Item(116, 9).Value = 21672776  # I116: 20839346 -> 21672776
$ws.Cells.Item(116, 10).Value = 166673330  # J116: 125005750 -> 166673330
$ws.Cells.Item(116, 11).Value = 21672776  # K116: 20839346 -> 21672776
$ws.Cells.Item(116, 12).Value = 166673330  # L116: 125005750 -> 166673330
$ws.Cells.Item(116, 13).Value = -21669334  # M116: -20835904 -> -21669334
$ws.Cells.Item(116, 14).Value = -166680214  # N116: -125012634 -> -166680214
# Row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Cells.Item(132, 8).Value = 13862  # H132: 14561.948 -> 13862
$ws.Cells.Item(132, 9).Value = 6056.4  # I132: 6475.107 -> 6056.4
$ws.Cells.Item(132, 10).Value = 21415.807  # J132: 22109.666 -> 21415.807
$ws.Cells.Item(132, 11).Value = 18169.2  # K132: 19425.321 -> 18169.2
$ws.Cells.Item(132, 12).Value = 64247.421  # L132: 66328.99800000001 -> 64247.421
$ws.Cells.Item(132, 13).Value = -15639.2  # M132: -16895.321 -> -15639.2
$ws.Cells.Item(132, 14).Value = -69307.421  # N132: -71388.99800000001 -> -69307.421
# Row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws.Cells.Item(137, 8).Value = 14932105  # H137: 13737897 -> 14932105
$ws.Cells.Item(137, 9).Value = 1250637.6  # I137: 1001410.1 -> 1250637.6
$ws.Cells.Item(137, 11).Value = 3751912.8  # K137: 3004230.3 -> 3751912.8
$ws.Cells.Item(137, 13).Value = -3749362.8  # M137: -3001680.3 -> -3749362.8

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots / Bronze Ingot
$ws.Cells.Item(2, 8).Value = 833821.4399999999  # H2: 875463 -> 833821.4399999999
$ws.Cells.Item(2, 9).Value = 1029497.4  # I2: 1093779 -> 1029497.4
$ws.Cells.Item(2, 11).Value = 1029497.4  # K2: 1093779 -> 1029497.4
$ws.Cells.Item(2, 13).Value = -1029384.4  # M2: -1093666 -> -1029384.4
# Row 32: Ingot We Trust / Steel Ingot
$ws.Cells.Item(32, 8).Value = 11659.213  # H32: 11668.487 -> 11659.213
$ws.Cells.Item(32, 9).Value = 11967.984  # I32: 11979.578 -> 11967.984
$ws.Cells.Item(32, 11).Value = 11967.984  # K32: 11979.578 -> 11967.984
$ws.Cells.Item(32, 13).Value = -11680.984  # M32: -11692.578 -> -11680.984
# Row 59: Parasitic Win / Cobalt-plated Caligae
$ws.Cells.Item(59, 8).Value = 0  # H59: 15000 -> 0
$ws.Cells.Item(59, 10).Value = 0  # J59: 15000 -> 0
$ws.Cells.Item(59, 12).Value = 0  # L59: 15000 -> 0
$ws.Cells.Item(59, 14).ClearContents()  # N59: -16608 -> (cleared)
# Row 60: Booty Call / Cobalt-plated Jackboots
$ws.Cells.Item(60, 8).Value = 15000  # H60: 0 -> 15000
$ws.Cells.Item(60, 10).Value = 15000  # J60: 0 -> 15000
$ws.Cells.Item(60, 12).Value = 15000  # L60: 0 -> 15000
$ws.Cells.Item(60, 14).Value = -16466  # N60: None -> -16466
# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Cells.Item(61, 8).Value = 6525.6904  # H61: 6731.39 -> 6525.6904
$ws.Cells.Item(61, 9).Value = 7612.6562  # I61: 7919.7744 -> 7612.6562
$ws.Cells.Item(61, 11).Value = 7612.6562  # K61: 7919.7744 -> 7612.6562
$ws.Cells.Item(61, 13).Value = -7400.6562  # M61: -7707.7744 -> -7400.6562
# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Cells.Item(74, 8).Value = 6098737.5  # H74: 7577054.5 -> 6098737.5
$ws.Cells.Item(74, 9).Value = 10000910  # I74: 13158902 -> 10000910
$ws.Cells.Item(74, 10).Value = 1592.6875  # J74: 1689.8572 -> 1592.6875
$ws.Cells.Item(74, 11).Value = 10000910  # K74: 13158902 -> 10000910
$ws.Cells.Item(74, 12).Value = 1592.6875  # L74: 1689.8572 -> 1592.6875
$ws.Cells.Item(74, 13).Value = -10000036  # M74: -13158028 -> -10000036
$ws.Cells.Item(74, 14).Value = -3340.6875  # N74: -3437.8572 -> -3340.6875
# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Cells.Item(77, 8).Value = 6098737.5  # H77: 7577054.5 -> 6098737.5
$ws.Cells.Item(77, 9).Value = 10000910  # I77: 13158902 -> 10000910
$ws.Cells.Item(77, 10).Value = 1592.6875  # J77: 1689.8572 -> 1592.6875
$ws.Cells.Item(77, 11).Value = 50004550  # K77: 65794510 -> 50004550
$ws.Cells.Item(77, 12).Value = 7963.4375  # L77: 8449.286 -> 7963.4375
$ws.Cells.Item(77, 13).Value = -50000182  # M77: -65790142 -> -50000182
$ws.Cells.Item(77, 14).Value = -16699.4375  # N77: -17185.286 -> -16699.4375
# Row 110: Scheduled Maintenance / Deepgold Ingot
$ws.Cells.Item(110, 8).Value = 5110543.5  # H110: 6812387.5 -> 5110543.5
$ws.Cells.Item(110, 9).Value = 10206588  # I110: 20408164 -> 10206588
$ws.Cells.Item(110, 11).Value = 10206588  # K110: 20408164 -> 10206588
$ws.Cells.Item(110, 13).Value = -10204543  # M110: -20406119 -> -10204543
# Row 116: No Scope / Titanbronze Ingot
$ws.Cells.Item(116, 8).Value = 833821.4399999999  # H116: 875463 -> 833821.4399999999
$ws.Cells.Item(116, 9).Value = 1029497.4  # I116: 1093779 -> 1029497.4
$ws.Cells.Item(116, 11).Value = 1029497.4  # K116: 1093779 -> 1029497.4
$ws.Cells.Item(116, 13).Value = -1027203.4  # M116: -1091485 -> -1027203.4
# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Cells.Item(136, 8).Value = 6525.6904  # H136: 6731.39 -> 6525.6904
$ws.Cells.Item(136, 9).Value = 7612.6562  # I136: 7919.7744 -> 7612.6562
$ws.Cells.Item(136, 11).Value = 22837.9686  # K136: 23759.3232 -> 22837.9686
$ws.Cells.Item(136, 13).Value = -20287.9686  # M136: -21209.3232 -> -20287.9686

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells / Bronze Ingot
$ws.Cells.Item(3, 8).Value = 833821.4399999999  # H3: 875463 -> 833821.4399999999
$ws.Cells.Item(3, 9).Value = 1029497.4  # I3: 1093779 -> 1029497.4
$ws.Cells.Item(3, 11).Value = 1029497.4  # K3: 1093779 -> 1029497.4
$ws.Cells.Item(3, 13).Value = -1029383.4  # M3: -1093665 -> -1029383.4
# Row 20: Smelt and Dealt / Iron Ingot
$ws.Cells.Item(20, 8).Value = 7652  # H20: 8691 -> 7652
$ws.Cells.Item(20, 9).Value = 8418  # I20: 8816.556 -> 8418
$ws.Cells.Item(20, 10).Value = 6790.25  # J20: 8529.571 -> 6790.25
$ws.Cells.Item(20, 11).Value = 8418  # K20: 8816.556 -> 8418
$ws.Cells.Item(20, 12).Value = 6790.25  # L20: 8529.571 -> 6790.25
$ws.Cells.Item(20, 13).Value = -8171  # M20: -8569.556 -> -8171
$ws.Cells.Item(20, 14).Value = -7284.25  # N20: -9023.571 -> -7284.25
# Row 107: The Gold Experience / Deepgold Nugget
$ws.Cells.Item(107, 8).Value = 7250  # H107: 7199.8 -> 7250
$ws.Cells.Item(107, 9).Value = 7500  # I107: 7333 -> 7500
$ws.Cells.Item(107, 11).Value = 7500  # K107: 7333 -> 7500
$ws.Cells.Item(107, 13).Value = -5580  # M107: -5413 -> -5580
# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Cells.Item(134, 8).Value = 1351.091  # H134: 1376.1364 -> 1351.091
$ws.Cells.Item(134, 9).Value = 983.45  # I134: 1011 -> 983.45
$ws.Cells.Item(134, 11).Value = 2950.35  # K134: 3033 -> 2950.35
$ws.Cells.Item(134, 13).Value = -415.3500000000004  # M134: -498 -> -415.3500000000004

$ws = $wb.Worksheets.Item("CRP")
# Row 16: Raise the Roof / Ash Lumber
$ws.Cells.Item(16, 8).Value = 2752.75  # H16: 2326.2222 -> 2752.75
$ws.Cells.Item(16, 9).Value = 2755.5  # I16: 2205.1428 -> 2755.5
$ws.Cells.Item(16, 11).Value = 2755.5  # K16: 2205.1428 -> 2755.5
$ws.Cells.Item(16, 13).Value = -2468.5  # M16: -1918.1428 -> -2468.5
# Row 31: Wall Not Found / Walnut Lumber
$ws.Cells.Item(31, 8).Value = 5394.5806  # H31: 5502.918 -> 5394.5806
$ws.Cells.Item(31, 9).Value = 1382.4445  # I31: 1480.4445 -> 1382.4445
$ws.Cells.Item(31, 10).Value = 7035.909  # J31: 7186.744 -> 7035.909
$ws.Cells.Item(31, 11).Value = 1382.4445  # K31: 1480.4445 -> 1382.4445
$ws.Cells.Item(31, 12).Value = 7035.909  # L31: 7186.744 -> 7035.909
$ws.Cells.Item(31, 13).Value = -1087.4445  # M31: -1185.4445 -> -1087.4445
$ws.Cells.Item(31, 14).Value = -7625.909  # N31: -7776.744 -> -7625.909
# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Cells.Item(34, 8).Value = 5394.5806  # H34: 5502.918 -> 5394.5806
$ws.Cells.Item(34, 9).Value = 1382.4445  # I34: 1480.4445 -> 1382.4445
$ws.Cells.Item(34, 10).Value = 7035.909  # J34: 7186.744 -> 7035.909
$ws.Cells.Item(34, 11).Value = 1382.4445  # K34: 1480.4445 -> 1382.4445
$ws.Cells.Item(34, 12).Value = 7035.909  # L34: 7186.744 -> 7035.909
$ws.Cells.Item(34, 13).Value = -1180.4445  # M34: -1278.4445 -> -1180.4445
$ws.Cells.Item(34, 14).Value = -7439.909  # N34: -7590.744 -> -7439.909
# Row 99: O Pine / Pine Lumber
$ws.Cells.Item(99, 8).Value = 13961.454  # H99: 11854.363 -> 13961.454
$ws.Cells.Item(99, 9).Value = 13961.454  # I99: 15059.9 -> 13961.454
$ws.Cells.Item(99, 10).Value = 0  # J99: 9183.083000000001 -> 0
$ws.Cells.Item(99, 11).Value = 13961.454  # K99: 15059.9 -> 13961.454
$ws.Cells.Item(99, 12).Value = 0  # L99: 9183.083000000001 -> 0
$ws.Cells.Item(99, 13).Value = -12463.454  # M99: -13561.9 -> -12463.454
$ws.Cells.Item(99, 14).ClearContents()  # N99: -12179.083 -> (cleared)
# Row 113: Patient Patients / White Ash Lumber
$ws.Cells.Item(113, 8).Value = 2752.75  # H113: 2326.2222 -> 2752.75
$ws.Cells.Item(113, 9).Value = 2755.5  # I113: 2205.1428 -> 2755.5
$ws.Cells.Item(113, 11).Value = 2755.5  # K113: 2205.1428 -> 2755.5
$ws.Cells.Item(113, 13).Value = -585.5  # M113: -35.14280000000008 -> -585.5
# Row 126: A Better Conductor / Red Pine Lumber
$ws.Cells.Item(126, 8).Value = 13961.454  # H126: 11854.363 -> 13961.454
$ws.Cells.Item(126, 9).Value = 13961.454  # I126: 15059.9 -> 13961.454
$ws.Cells.Item(126, 10).Value = 0  # J126: 9183.083000000001 -> 0
$ws.Cells.Item(126, 11).Value = 41884.362  # K126: 45179.7 -> 41884.362
$ws.Cells.Item(126, 12).Value = 0  # L126: 27549.249 -> 0
$ws.Cells.Item(126, 13).Value = -39414.362  # M126: -42709.7 -> -39414.362
$ws.Cells.Item(126, 14).ClearContents()  # N126: -32489.249 -> (cleared)
# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Cells.Item(132, 8).Value = 27808440  # H132: 37075530 -> 27808440
$ws.Cells.Item(132, 9).Value = 33359358  # I132: 47653150 -> 33359358
$ws.Cells.Item(132, 11).Value = 100078074  # K132: 142959450 -> 100078074
$ws.Cells.Item(132, 13).Value = -100075544  # M132: -142956920 -> -100075544

$ws = $wb.Worksheets.Item("CUL")
# Row 88: Don't Let It Fall Apart / Liver-cheese Sandwich
$ws.Cells.Item(88, 8).Value = 10750  # H88: 9197.4 -> 10750
$ws.Cells.Item(88, 10).Value = 10750  # J88: 9197.4 -> 10750
$ws.Cells.Item(88, 12).Value = 32250  # L88: 27592.2 -> 32250
$ws.Cells.Item(88, 14).Value = -33106  # N88: -28448.2 -> -33106
# Row 91: Better Come Back with a Sandwich (L) / Liver-cheese Sandwich
$ws.Cells.Item(91, 8).Value = 10750  # H91: 9197.4 -> 10750
$ws.Cells.Item(91, 10).Value = 10750  # J91: 9197.4 -> 10750
$ws.Cells.Item(91, 12).Value = 32250  # L91: 27592.2 -> 32250
$ws.Cells.Item(91, 14).Value = -35214  # N91: -30556.2 -> -35214
# Row 120: A Happy End / Paella
$ws.Cells.Item(120, 8).Value = 23756.25  # H120: 26575.715 -> 23756.25
$ws.Cells.Item(120, 9).Value = 5016.6665  # I120: 5515 -> 5016.6665
$ws.Cells.Item(120, 11).Value = 15049.9995  # K120: 16545 -> 15049.9995
$ws.Cells.Item(120, 13).Value = -10211.9995  # M120: -11707 -> -10211.9995
# Row 134: Don't Knock It Till You've Tried It / Mezcal-marinated Swampmonk
$ws.Cells.Item(134, 8).Value = 999  # H134: 10483.583 -> 999
$ws.Cells.Item(134, 9).Value = 0  # I134: 3336.2222 -> 0
$ws.Cells.Item(134, 10).Value = 999  # J134: 31925.666 -> 999
$ws.Cells.Item(134, 11).Value = 0  # K134: 10008.6666 -> 0
$ws.Cells.Item(134, 12).Value = 2997  # L134: 95776.99800000001 -> 2997
$ws.Cells.Item(134, 13).ClearContents()  # M134: -4938.6666 -> (cleared)
$ws.Cells.Item(134, 14).Value = -13137  # N134: -105916.998 -> -13137

$ws = $wb.Worksheets.Item("GSM")
# Row 59: Sew Not Doing This / Electrum Needle
$ws.Cells.Item(59, 8).Value = 48999  # H59: 49000 -> 48999
$ws.Cells.Item(59, 10).Value = 48999  # J59: 49000 -> 48999
$ws.Cells.Item(59, 12).Value = 48999  # L59: 49000 -> 48999
$ws.Cells.Item(59, 14).Value = -50165  # N59: -50166 -> -50165
# Row 102: Put the Metal to the Peddle / Durium Ingot
$ws.Cells.Item(102, 8).Value = 23818624  # H102: 25009436 -> 23818624
$ws.Cells.Item(102, 9).Value = 31258944  # I102: 33342716 -> 31258944
$ws.Cells.Item(102, 11).Value = 31258944  # K102: 33342716 -> 31258944
$ws.Cells.Item(102, 13).Value = -31257322  # M102: -33341094 -> -31257322

$ws = $wb.Worksheets.Item("LTW")
# Row 122: Hell on Leather / Gaja Leather
$ws.Cells.Item(122, 8).Value = 47624836  # H122: 50005830 -> 47624836
$ws.Cells.Item(122, 9).Value = 76927656  # I122: 83337890 -> 76927656
$ws.Cells.Item(122, 11).Value = 230782968  # K122: 250013670 -> 230782968
$ws.Cells.Item(122, 13).Value = -230780518  # M122: -250011220 -> -230780518
# Row 136: Respect for Br'aax / Br'aax Leather
$ws.Cells.Item(136, 8).Value = 3139.3684  # H136: 2884.3125 -> 3139.3684
$ws.Cells.Item(136, 9).Value = 2663.2  # I136: 2173 -> 2663.2
$ws.Cells.Item(136, 10).Value = 4925  # J136: 5966.6665 -> 4925
$ws.Cells.Item(136, 11).Value = 7989.599999999999  # K136: 6519 -> 7989.599999999999
$ws.Cells.Item(136, 12).Value = 14775  # L136: 17899.9995 -> 14775
$ws.Cells.Item(136, 13).Value = -5439.599999999999  # M136: -3969 -> -5439.599999999999
$ws.Cells.Item(136, 14).Value = -19875  # N136: -22999.9995 -> -19875

$ws = $wb.Worksheets.Item("WVR")
# Row 96: Skills on Display / Ruby Cotton Cloth
$ws.Cells.Item(96, 8).Value = 5039  # H96: 5039.8887 -> 5039
$ws.Cells.Item(96, 10).Value = 4915.4  # J96: 4917 -> 4915.4
$ws.Cells.Item(96, 12).Value = 4915.4  # L96: 4917 -> 4915.4
$ws.Cells.Item(96, 14).Value = -7661.4  # N96: -7663 -> -7661.4
# Row 126: A Polished Purchase / Snow Linen
$ws.Cells.Item(126, 8).Value = 4397.9  # H126: 4557.4 -> 4397.9
$ws.Cells.Item(126, 9).Value = 3800  # I126: 4581.5 -> 3800
$ws.Cells.Item(126, 10).Value = 4547.375  # J126: 4541.3335 -> 4547.375
$ws.Cells.Item(126, 11).Value = 11400  # K126: 13744.5 -> 11400
$ws.Cells.Item(126, 12).Value = 13642.125  # L126: 13624.0005 -> 13642.125
$ws.Cells.Item(126, 13).Value = -8930  # M126: -11274.5 -> -8930
$ws.Cells.Item(126, 14).Value = -18582.125  # N126: -18564.0005 -> -18582.125
# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Cells.Item(132, 8).Value = 10207366  # H132: 10641701 -> 10207366
$ws.Cells.Item(132, 9).Value = 1144.1613  # I132: 1165.6333 -> 1144.1613
$ws.Cells.Item(132, 10).Value = 27784750  # J132: 29419116 -> 27784750
$ws.Cells.Item(132, 11).Value = 3432.4839  # K132: 3496.8999 -> 3432.4839
$ws.Cells.Item(132, 12).Value = 83354250  # L132: 88257348 -> 83354250
$ws.Cells.Item(132, 13).Value = -902.4839000000002  # M132: -966.8998999999999 -> -902.4839000000002
$ws.Cells.Item(132, 14).Value = -83359310  # N132: -88262408 -> -83359310
# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Cells.Item(136, 8).Value = 7957.458  # H136: 8249.1875 -> 7957.458
$ws.Cells.Item(136, 9).Value = 2441.9285  # I136: 2713.6 -> 2441.9285
$ws.Cells.Item(136, 11).Value = 7325.7855  # K136: 8140.799999999999 -> 7325.7855
$ws.Cells.Item(136, 13).Value = -4775.7855  # M136: -5590.799999999999 -> -4775.7855
